# Update cryptocurrency price/volume data per upstream GitHub Actions scrape refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "59.746.49"
$ws.Range("E2").Value = "  -0.70%  "
$ws.Range("D3").Value = "2.305.46"
$ws.Range("E3").Value = "  -1.83%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.10%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "542.78"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.93%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "128.56"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -4.07%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.999"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.14%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.570"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -3.52%  "
$ws.Range("D9").Value = "2.304.49"
$ws.Range("E9").Value = "  -1.73%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.101"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.44%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.54"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.21%  "
$ws.Range("E13").Value = "  -1.96%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "23.15"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -4.33%  "
$ws.Range("D15").Value = "2.714.64"
$ws.Range("E15").Value = "  -1.91%  "
$ws.Range("D16").Value = "59.561.43"
$ws.Range("E16").Value = "  -0.73%  "
$ws.Range("E17").Value = "  -2.11%  "
$ws.Range("E18").Value = "  -2.12%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "10.41"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -3.03%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "4.03"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -4.51%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "310.23"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -2.24%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.51"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -4.64%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.999"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.64%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "63.01"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.25%  "
$ws.Range("E25").Value = "  -3.33%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.999"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.10%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "7.73"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -4.21%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.35"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.99%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.20"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +3.03%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "171.71"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.09%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.71"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -2.32%  "
$ws.Range("D32").Value = "0.0₃0715"
$ws.Range("E32").Value = "  -5.41%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.80"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -2.34%  "
$ws.Range("E34").Value = "  -3.30%  "
$ws.Range("E35").Value = "  +0.00%  "
$ws.Range("E36").Value = "  -7.60%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "17.62"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -2.24%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.00"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.07%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.99"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -3.40%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "313.89"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -1.26%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "37.48"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -2.02%  "
$ws.Range("E42").Value = "  -3.41%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "136.15"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -6.37%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "3.42"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -1.77%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0939"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -2.24%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.567"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.62%  "
$ws.Range("B47").Value = "InjectiveProtocol"
$ws.Range("C47").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "18.58"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -1.40%  "
$ws.Range("B48").Value = "BabyDogeCoin"
$ws.Range("C48").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D48").Value = "0.0₆0229"
$ws.Range("E48").Value = "  +27.62%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0488"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -2.30%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0212"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.58%  "
$ws.Range("E51").Value = "  -0.32%  "
